$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "-"

$ws.Range("B6").Value = "MCT-1A-Gestão integrada"
$ws.Range("E6").Value = "MEC-1B-Gestao Intregrada"

$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "MEC-1B-Gestao Intregrada"

$ws.Range("D11").Value = "MEC-1A-Gestao Integrada"
$ws.Range("D12").Value = "MEC-1A-Gestao Integrada"

$ws.Range("B14").Value = "-"
$ws.Range("E14").Value = "-"
